# Update rows 2-4 (ECs sending cluster) with refreshed TPM-derived values
# and append rows 5-7 (MuSCs sending cluster), per the new sharedStrings
# ordering (ECs, MuSCs, Efna3, Epha7, FAPs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna3"
$ws.Cells.Item(2, 3).Value = "Epha7"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.046374
$ws.Cells.Item(2, 8).Value = 0.139122
$ws.Cells.Item(2, 9).Value = 0.6592866045237633
$ws.Cells.Item(2, 10).Value = 0.6592866045237632
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.05661333333333334
$ws.Cells.Item(2, 14).Value = 0.16984
$ws.Cells.Item(2, 15).Value = 0.0204119846136133
$ws.Cells.Item(2, 16).Value = 0.02041198461361329
$ws.Cells.Item(2, 17).Value = 0.00262538672
$ws.Cells.Item(2, 18).Value = 0.02362848048
$ws.Cells.Item(2, 19).Value = 0.01345734802750041
$ws.Cells.Item(2, 20).Value = 0.01345734802750041

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna3"
$ws.Cells.Item(3, 3).Value = "Epha7"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.046374
$ws.Cells.Item(3, 8).Value = 0.139122
$ws.Cells.Item(3, 9).Value = 0.6592866045237633
$ws.Cells.Item(3, 10).Value = 0.6592866045237632
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.2780386666666667
$ws.Cells.Item(3, 14).Value = 0.8341160000000001
$ws.Cells.Item(3, 15).Value = 0.100247073468963
$ws.Cells.Item(3, 16).Value = 0.1002470734689629
$ws.Cells.Item(3, 17).Value = 0.012893765128
$ws.Cells.Item(3, 18).Value = 0.116043886152
$ws.Cells.Item(3, 19).Value = 0.06609155268079683
$ws.Cells.Item(3, 20).Value = 0.0660915526807968

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna3"
$ws.Cells.Item(4, 3).Value = "Epha7"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.046374
$ws.Cells.Item(4, 8).Value = 0.139122
$ws.Cells.Item(4, 9).Value = 0.6592866045237633
$ws.Cells.Item(4, 10).Value = 0.6592866045237632
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.438882
$ws.Cells.Item(4, 14).Value = 7.316646
$ws.Cells.Item(4, 15).Value = 0.8793409419174237
$ws.Cells.Item(4, 16).Value = 0.8793409419174237
$ws.Cells.Item(4, 17).Value = 0.113100713868
$ws.Cells.Item(4, 18).Value = 1.017906424812
$ws.Cells.Item(4, 19).Value = 0.5797377038154661
$ws.Cells.Item(4, 20).Value = 0.579737703815466

# Row 5
$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Efna3"
$ws.Cells.Item(5, 3).Value = "Epha7"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.02396566666666667
$ws.Cells.Item(5, 8).Value = 0.071897
$ws.Cells.Item(5, 9).Value = 0.3407133954762367
$ws.Cells.Item(5, 10).Value = 0.3407133954762367
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.05661333333333334
$ws.Cells.Item(5, 14).Value = 0.16984
$ws.Cells.Item(5, 15).Value = 0.0204119846136133
$ws.Cells.Item(5, 16).Value = 0.02041198461361329
$ws.Cells.Item(5, 17).Value = 0.001356776275555556
$ws.Cells.Item(5, 18).Value = 0.01221098648
$ws.Cells.Item(5, 19).Value = 0.006954636586112887
$ws.Cells.Item(5, 20).Value = 0.006954636586112886

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Efna3"
$ws.Cells.Item(6, 3).Value = "Epha7"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.02396566666666667
$ws.Cells.Item(6, 8).Value = 0.071897
$ws.Cells.Item(6, 9).Value = 0.3407133954762367
$ws.Cells.Item(6, 10).Value = 0.3407133954762367
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.2780386666666667
$ws.Cells.Item(6, 14).Value = 0.8341160000000001
$ws.Cells.Item(6, 15).Value = 0.100247073468963
$ws.Cells.Item(6, 16).Value = 0.1002470734689629
$ws.Cells.Item(6, 17).Value = 0.006663382005777779
$ws.Cells.Item(6, 18).Value = 0.05997043805200001
$ws.Cells.Item(6, 19).Value = 0.03415552078816614
$ws.Cells.Item(6, 20).Value = 0.03415552078816613

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Efna3"
$ws.Cells.Item(7, 3).Value = "Epha7"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.02396566666666667
$ws.Cells.Item(7, 8).Value = 0.071897
$ws.Cells.Item(7, 9).Value = 0.3407133954762367
$ws.Cells.Item(7, 10).Value = 0.3407133954762367
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.438882
$ws.Cells.Item(7, 14).Value = 7.316646
$ws.Cells.Item(7, 15).Value = 0.8793409419174237
$ws.Cells.Item(7, 16).Value = 0.8793409419174237
$ws.Cells.Item(7, 17).Value = 0.05844943305133334
$ws.Cells.Item(7, 18).Value = 0.526044897462
$ws.Cells.Item(7, 19).Value = 0.2996032381019577
$ws.Cells.Item(7, 20).Value = 0.2996032381019577
